$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add Box Number values into column A for rows 2..56 ("01".."55"),
# matching the same text style already used by column C (numFmtId 49 -> "@" text format).
for ($row = 2; $row -le 56; $row++) {
    $boxNumber = "{0:D2}" -f ($row - 1)
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $boxNumber
}

# Update the active selection to I21, matching the saved sheet view state.
$ws.Range("I21").Select()
